# Fixed vocab URI in exemple 7
# Change the ConceptScheme URI from ".../days" to ".../paintings", turn it
# into a hyperlink (matching the existing hyperlinks on C2/C3), and touch up
# the surrounding formatting the same way the original author's edit did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the ConceptScheme URI value -----------------------------
$newUri = "http://data.sparna.fr/vocabularies/paintings"
$ws.Range("B1").Value = $newUri

# --- 2. Turn B1 into a hyperlink, like the PREFIX cells below it -------
$ws.Hyperlinks.Add($ws.Range("B1"), $newUri, "", "", $newUri) | Out-Null

# Adding a hyperlink normally stamps the built-in blue/underlined
# "Hyperlink" style onto the cell; the source file doesn't use that look
# for its other hyperlinks (C2/C3), so put B1's font back the way it was.
$ws.Range("B1").Font.Name = "Calibri"
$ws.Range("B1").Font.Size = 11
$ws.Range("B1").Font.Bold = $false
$ws.Range("B1").Font.Underline = 0
$ws.Range("B1").Font.Color = 0

# --- 3. Row 1 grew slightly taller after the edit -----------------------
$ws.Rows("1:1").RowHeight = 14.95

# --- 4. Columns picked up a small width bump as well --------------------
$ws.Columns("A:A").ColumnWidth = 40.514844804318464
$ws.Columns("B:B").ColumnWidth = 23.37719298245617
$ws.Columns("C:C").ColumnWidth = 31.838731443994565
$ws.Columns("D:D").ColumnWidth = 21.126180836707167
$ws.Columns("E:E").ColumnWidth = 35.373144399460166
$ws.Columns("F:F").ColumnWidth = 28.409581646423767

# --- 5. Column D picks up blank formatted cells on rows 5, 9 and 10 -----
$ws.Range("D1:D10").NumberFormat = "General"

# --- 6. Selection ends up parked on B5 -----------------------------------
$ws.Range("B5").Select() | Out-Null
